# spring 23 week 7 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.43
$ws.Range("E2").Value = 1.26
$ws.Range("F2").Value = 0.71

# Row 3
$ws.Range("B3").Value = 1.4

# Row 4
$ws.Range("C4").Value = 1.43
$ws.Range("E4").Value = 1.21

# Row 5
$ws.Range("B5").Value = 1.47
$ws.Range("D5").Value = 1.35
$ws.Range("G5").Value = 0.63

# Row 6
$ws.Range("B6").Value = 2.07
$ws.Range("E6").Value = 1.29
$ws.Range("F6").Value = 1.14
$ws.Range("G6").Value = 1.01

# Row 7
$ws.Range("E7").Value = 1.98
$ws.Range("F7").Value = 1.53
